$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Timetable")

# 1. AG2: room C002 -> C003 for the HS161 English Language and Communication entry
$ws.Range("AG2").Value = "HS161 | English Language and Communication | Dr. Rajesh N S | C003"

# 2. Row 5 rework: the CS161 slot (was I5:N5, room C003) moves later and splits,
#    making room for a LUNCH BREAK block at N5:Q5 and the CS161 class (room C002) at R5:W5.

# First split the existing I5:N5 merge back into individual cells.
$ws.Range("I5:N5").UnMerge()

# Clear the old CS161 text out of I5 and blank out I5:M5 (style + content),
# using the already-blank H5 as a style donor so they end up with no explicit style.
$ws.Range("H5").Copy($ws.Range("I5:M5"))
$ws.Range("I5:M5").ClearContents()

# N5 becomes the new LUNCH BREAK cell - borrow the LUNCH BREAK style from O5 first
# (before we blank it out below), then set its text.
$ws.Range("O5").Copy($ws.Range("N5"))
$ws.Range("N5").Value = "LUNCH BREAK"

# Clear the other old LUNCH BREAK cells that are now absorbed into the N5:Q5 merge.
$ws.Range("O5:Q5").ClearContents()

# R5 becomes the CS161 class cell (room C002 now) - borrow the course-slot style from I2.
$ws.Range("I2").Copy($ws.Range("R5"))
$ws.Range("R5").Value = "CS161 | Problem Solving through Programming | Dr. Sunil C K | C002"

# Clear the remaining old LUNCH BREAK cells that are now absorbed into the R5:W5 merge.
$ws.Range("S5:W5").ClearContents()

# Recreate the merges: N5:Q5 (LUNCH BREAK) and R5:W5 (CS161 class).
$ws.Range("N5:Q5").Merge()
$ws.Range("R5:W5").Merge()

# Merging (like the real Excel UI) paints the top-left cell's style across the
# whole merged range; strip that residual formatting from the non-anchor cells
# so only the anchor cell (N5 / R5) carries the style, matching the rest of
# the sheet's merged ranges (e.g. B2:G2).
$ws.Range("O5:Q5").ClearFormats()
$ws.Range("S5:W5").ClearFormats()
